$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.774353504180908
$ws.Range("F3").Value = 0.07007050514221191
$ws.Range("F4").Value = 0.001367568969726562
$ws.Range("F5").Value = 0.001143693923950195
$ws.Range("F6").Value = 0.03132462501525879
